$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Row 141 ---
$ws.Range("B141").Value = "SingleUseId268"
$ws.Range("C141").Value = "Default"
$ws.Range("D141").Value = "Left"
$ws.Range("E141").Value = "LTR"
$ws.Range("F141").Value = "<value>"

# --- Row 142 ---
$ws.Range("B142").Value = "SingleUseId269"
$ws.Range("C142").Value = "Default"
$ws.Range("D142").Value = "Left"
$ws.Range("E142").Value = "LTR"
$ws.Range("F142").Value = "<value>"

# --- Row 143 ---
$ws.Range("B143").Value = "SingleUseId270"
$ws.Range("C143").Value = "Default"
$ws.Range("D143").Value = "Left"
$ws.Range("E143").Value = "LTR"

# "1000" must stay a text value (not be coerced to a number), so type it
# into a scratch cell formatted as Text, copy it, and paste-values it into
# place - this carries the text type over without leaving the Text number
# format behind on the destination cell.
$scratch = $ws.Range("ZZ500")
$scratch.NumberFormat = "@"
$scratch.Value = "1000"
$scratch.Copy()
$ws.Range("F143").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()

# --- Row 144 ---
$ws.Range("B144").Value = "SingleUseId271"
$ws.Range("C144").Value = "Default"
$ws.Range("D144").Value = "Left"
$ws.Range("E144").Value = "LTR"

$scratch = $ws.Range("ZZ500")
$scratch.NumberFormat = "@"
$scratch.Value = "-1000"
$scratch.Copy()
$ws.Range("F144").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$scratch.Clear()
